# daily auto push: 2026-02-02 10:03 UTC
# A new daily record (2026/02/02, 月, 16, 25) was inserted into the sorted
# log at what becomes row 743, pushing every following row down by one and
# appending a fresh final row (785) that carries the data which used to
# belong to the old last row (784).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 743..784 down to 744..785, opening up a blank row 743.
$ws.Rows.Item(743).Insert()

# Fill the newly opened row. The date column holds plain text like all the
# other rows ("2026/12/29", etc.) rather than a real date, so force a text
# entry (leading apostrophe) to stop Excel's automatic date detection, then
# strip the stray number-format Excel attaches when it believes it saw a
# date so the cell ends up styled exactly like its neighbours.
$ws.Range("A743").Value = "'2026/02/02"
$ws.Range("A743").ClearFormats()
$ws.Range("B743").Value = "月"
$ws.Range("C743").Value = 16
$ws.Range("D743").Value = 25
